# Apply the updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.825.55'
$ws.Range('E2').Value = '  +0.63%  '
# Row 3
$ws.Range('D3').Value = '3.318.23'
$ws.Range('E3').Value = '  +2.48%  '
# Row 4
$ws.Range('E4').Value = '  -0.01%  '
# Row 5
$ws.Range('D5').Value = "'605.41"
$ws.Range('E5').Value = '  +1.94%  '
# Row 6
$ws.Range('D6').Value = "'142.73"
$ws.Range('E6').Value = '  +0.70%  '
# Row 7
$ws.Range('E7').Value = '  -0.03%  '
# Row 8
$ws.Range('D8').Value = '3.317.19'
$ws.Range('E8').Value = '  +2.56%  '
# Row 9
$ws.Range('E9').Value = '  +0.02%  '
# Row 10
$ws.Range('E10').Value = '  +1.54%  '
# Row 11
$ws.Range('E11').Value = '  +4.07%  '
# Row 12
$ws.Range('E12').Value = '  +1.01%  '
# Row 13
$ws.Range('E13').Value = '  +0.50%  '
# Row 14
$ws.Range('E14').Value = '  +2.04%  '
# Row 15
$ws.Range('D15').Value = '3.866.61'
$ws.Range('E15').Value = '  +2.59%  '
# Row 16
$ws.Range('E16').Value = '  +0.27%  '
# Row 17
$ws.Range('D17').Value = '3.320.41'
$ws.Range('E17').Value = '  +2.62%  '
# Row 18
$ws.Range('D18').Value = '63.918.75'
$ws.Range('E18').Value = '  +0.84%  '
# Row 19
$ws.Range('D19').Value = "'6.87"
$ws.Range('E19').Value = '  +1.26%  '
# Row 20
$ws.Range('D20').Value = "'481.67"
$ws.Range('E20').Value = '  +1.25%  '
# Row 21
$ws.Range('D21').Value = "'14.11"
$ws.Range('E21').Value = '  -0.16%  '
# Row 22
$ws.Range('E22').Value = '  +1.81%  '
# Row 23
$ws.Range('E23').Value = '  +1.06%  '
# Row 24
$ws.Range('D24').Value = "'13.95"
$ws.Range('E24').Value = '  +5.88%  '
# Row 25
$ws.Range('D25').Value = "'85.06"
$ws.Range('E25').Value = '  +1.29%  '
# Row 26
$ws.Range('E26').Value = '  +0.09%  '
# Row 27
$ws.Range('E27').Value = '  +1.75%  '
# Row 28
$ws.Range('E28').Value = '  -0.06%  '
# Row 29
$ws.Range('D29').Value = "'8.26"
$ws.Range('E29').Value = '  +2.39%  '
# Row 30
$ws.Range('D30').Value = "'7.17"
$ws.Range('E30').Value = '  -4.28%  '
# Row 31
$ws.Range('D31').Value = "'2.15"
$ws.Range('E31').Value = '  +2.38%  '
# Row 32
$ws.Range('D32').Value = "'28.90"
$ws.Range('E32').Value = '  +5.08%  '
# Row 33
$ws.Range('D33').Value = "'0.107"
$ws.Range('E33').Value = '  -1.27%  '
# Row 34
$ws.Range('E34').Value = '  +0.10%  '
# Row 35
$ws.Range('E35').Value = '  +1.47%  '
# Row 36
$ws.Range('E36').Value = '  +3.13%  '
# Row 37
$ws.Range('D37').Value = "'52.44"
$ws.Range('E37').Value = '  -0.61%  '
# Row 38
$ws.Range('E38').Value = '  +4.53%  '
# Row 39
$ws.Range('E39').Value = '  +1.65%  '
# Row 40
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = "'434.41"
$ws.Range('E40').Value = '  +3.00%  '
# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.127.79'
$ws.Range('E41').Value = '  +5.03%  '
# Row 42
$ws.Range('E42').Value = '  +6.92%  '
# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = "'2.76"
$ws.Range('E43').Value = '  +0.11%  '
# Row 44
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D44').Value = "'8.37"
$ws.Range('E44').Value = '  -0.21%  '
# Row 45
$ws.Range('E45').Value = '  +0.29%  '
# Row 46
$ws.Range('E46').Value = '  +4.07%  '
# Row 47
$ws.Range('D47').Value = "'36.74"
$ws.Range('E47').Value = '  +7.23%  '
# Row 48
$ws.Range('D48').Value = "'26.39"
$ws.Range('E48').Value = '  +1.87%  '
# Row 50
$ws.Range('E50').Value = '  -1.07%  '
# Row 51
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = "'124.54"
$ws.Range('E51').Value = '  +2.61%  '
